$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44798
$ws.Range("J2").Value = 200

# Row 3
$ws.Range("D3").Value = 44847
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 7500
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 7750
$ws.Range("N3").Value = '$/docena de atados'
$ws.Range("P3").Value = 2583
$ws.Range("Q3").Value = 3

# Row 4
$ws.Range("D4").Value = 44782
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("P4").Value = 2667

# Row 5
$ws.Range("D5").Value = 44839
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 7500
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 7750
$ws.Range("N5").Value = '$/docena de atados'
$ws.Range("O5").Value = 'Provincia de Diguillín'
$ws.Range("P5").Value = 2583
$ws.Range("Q5").Value = 3

# Row 6
$ws.Range("D6").Value = 44764
$ws.Range("J6").Value = 100
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 8500
$ws.Range("P6").Value = 2833

# Row 7
$ws.Range("D7").Value = 44662

# Row 8
$ws.Range("D8").Value = 44804
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 8500
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8750
$ws.Range("P8").Value = 2917

# Row 9
$ws.Range("D9").Value = 44215
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 3500
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = 3768
$ws.Range("N9").Value = '$/paquete 2 kilos'
$ws.Range("O9").Value = 'Provincia de Diguillín'
$ws.Range("P9").Value = 1884
$ws.Range("Q9").Value = 2

# Row 10
$ws.Range("D10").Value = 44790
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 8500
$ws.Range("M10").Value = 8250
$ws.Range("N10").Value = '$/docena de atados'
$ws.Range("P10").Value = 2750
$ws.Range("Q10").Value = 3

# Row 11
$ws.Range("D11").Value = 44841
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 7500
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = 7750
$ws.Range("P11").Value = 2583

# Row 12
$ws.Range("D12").Value = 44830
$ws.Range("K12").Value = 7500
$ws.Range("M12").Value = 7750
$ws.Range("O12").Value = 'Provincia de Diguillín'
$ws.Range("P12").Value = 2583

# Row 13
$ws.Range("D13").Value = 44811
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 8000
$ws.Range("M13").Value = 8500
$ws.Range("P13").Value = 2833

# Row 14
$ws.Range("D14").Value = 44762
$ws.Range("J14").Value = 60
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 8000
$ws.Range("N14").Value = '$/docena de atados'
$ws.Range("P14").Value = 2667
$ws.Range("Q14").Value = 3

# Row 15
$ws.Range("D15").Value = 44208
$ws.Range("J15").Value = 85
$ws.Range("K15").Value = 3700
$ws.Range("L15").Value = 4000
$ws.Range("M15").Value = 3824
$ws.Range("N15").Value = '$/paquete 2 kilos'
$ws.Range("P15").Value = 1912
$ws.Range("Q15").Value = 2

# Row 16
$ws.Range("D16").Value = 44704
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 6500
$ws.Range("M16").Value = 6250
$ws.Range("P16").Value = 174

# Row 17
$ws.Range("D17").Value = 44664
$ws.Range("J17").Value = 200
$ws.Range("L17").Value = 8500
$ws.Range("M17").Value = 8250
$ws.Range("N17").Value = '$/paquete 36 unidades'
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 229
$ws.Range("Q17").Value = 36

# Row 18
$ws.Range("D18").Value = 44210
$ws.Range("J18").Value = 105
$ws.Range("K18").Value = 3500
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = 3714
$ws.Range("N18").Value = '$/paquete 2 kilos'
$ws.Range("P18").Value = 1857
$ws.Range("Q18").Value = 2

# Row 19
$ws.Range("D19").Value = 44791
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 8500
$ws.Range("M19").Value = 8250
$ws.Range("N19").Value = '$/docena de atados'
$ws.Range("P19").Value = 2750
$ws.Range("Q19").Value = 3

# Row 20
$ws.Range("D20").Value = 44701
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 7500
$ws.Range("M20").Value = 7250
$ws.Range("N20").Value = '$/paquete 36 unidades'
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 201
$ws.Range("Q20").Value = 36

# Row 21
$ws.Range("D21").Value = 44160
$ws.Range("J21").Value = 43
$ws.Range("K21").Value = 3500
$ws.Range("L21").Value = 4000
$ws.Range("M21").Value = 3709
$ws.Range("N21").Value = '$/paquete 36 unidades'
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 103
$ws.Range("Q21").Value = 36

# Row 22
$ws.Range("D22").Value = 44784
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8500
$ws.Range("M22").Value = 8250
$ws.Range("N22").Value = '$/docena de atados'
$ws.Range("O22").Value = 'Provincia de Diguillín'
$ws.Range("P22").Value = 2750
$ws.Range("Q22").Value = 3

# Row 23
$ws.Range("D23").Value = 44161
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 2800
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = 2900
$ws.Range("N23").Value = '$/paquete 2 kilos'
$ws.Range("P23").Value = 1450
$ws.Range("Q23").Value = 2

# Row 24
$ws.Range("D24").Value = 44817
$ws.Range("J24").Value = 120
$ws.Range("L24").Value = 8500
$ws.Range("M24").Value = 8250
$ws.Range("P24").Value = 2750

# Row 25
$ws.Range("D25").Value = 44223
$ws.Range("J25").Value = 80
$ws.Range("L25").Value = 3800
$ws.Range("M25").Value = 3688
$ws.Range("P25").Value = 1844

# Row 26
$ws.Range("D26").Value = 44818
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 8500
$ws.Range("M26").Value = 8250
$ws.Range("N26").Value = '$/docena de atados'
$ws.Range("P26").Value = 2750
$ws.Range("Q26").Value = 3

# Row 27
$ws.Range("D27").Value = 44771
$ws.Range("J27").Value = 150
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 8000
$ws.Range("P27").Value = 2667

# Row 28
$ws.Range("D28").Value = 44810
$ws.Range("J28").Value = 120

# Row 29
$ws.Range("D29").Value = 44760
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 8000
$ws.Range("N29").Value = '$/docena de atados'
$ws.Range("P29").Value = 2667
$ws.Range("Q29").Value = 3

# Row 30
$ws.Range("D30").Value = 44769
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 8000
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 8000
$ws.Range("N30").Value = '$/docena de atados'
$ws.Range("P30").Value = 2667
$ws.Range("Q30").Value = 3

# Row 31
$ws.Range("D31").Value = 44166
$ws.Range("J31").Value = 70
$ws.Range("K31").Value = 3500
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = 3679
$ws.Range("N31").Value = '$/paquete 36 unidades'
$ws.Range("O31").Value = 'Región Metropolitana'
$ws.Range("P31").Value = 102
$ws.Range("Q31").Value = 36

# Row 32
$ws.Range("D32").Value = 44225
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 3400
$ws.Range("L32").Value = 3700
$ws.Range("M32").Value = 3550
$ws.Range("P32").Value = 1775

# Row 33
$ws.Range("D33").Value = 44775
$ws.Range("J33").Value = 100
$ws.Range("L33").Value = 8000
$ws.Range("M33").Value = 8000
$ws.Range("P33").Value = 2667

# Row 34
$ws.Range("D34").Value = 44209
$ws.Range("J34").Value = 150
$ws.Range("K34").Value = 3500
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = 3767
$ws.Range("N34").Value = '$/paquete 2 kilos'
$ws.Range("P34").Value = 1884
$ws.Range("Q34").Value = 2

# Row 35
$ws.Range("D35").Value = 44845
$ws.Range("K35").Value = 7500
$ws.Range("M35").Value = 7750
$ws.Range("P35").Value = 2583
